# Apply the commit "Update some slides and readmes" to the HTML.pptx deck:
#  1. Insert a new "Goal" slide (Title + bulleted Content placeholder) as the
#     new slide 2, pushing the existing "Creation of HTML" slide down to
#     slide 3.
#  2. Refresh the cached "datetimeFigureOut" footer field text on every
#     slide layout and the slide master from 7/19/2021 to 7/21/2021.

$p = $ppt.ActivePresentation

# --- 1. Insert the new "Goal" slide at position 2 ----------------------
# Layout 2 on the master is "Title and Content", the same layout used by
# the existing "Creation of HTML" slide.
$goal = $p.Slides.Add(2, 2)

$goal.Shapes.Item(1).TextFrame.TextRange.Text = "Goal"

$bullet1 = "We will not go over all the ins and outs of HTML as the possible elements and use cases are vast"
$bullet2 = "Go over some semantic HTML"
$bullet3 = "Discuss importance of using semantic HTML for SEO and accessibility"

$body = $goal.Shapes.Item(2).TextFrame.TextRange
$body.Text = $bullet1 + "`r" + $bullet2 + "`r" + $bullet3

# --- 2. Update the cached footer date text ------------------------------
function Update-DatePlaceholder($shapes) {
    foreach ($sh in $shapes) {
        if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "7/19/2021") {
            $sh.TextFrame.TextRange.Text = "7/21/2021"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $cl = $layouts.Item($i)
    Update-DatePlaceholder $cl.Shapes
}
